$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 600
$ws.Range("C3").Value = 1000
$ws.Range("C4").Value = 360
